# plantilla_datos_pacientes.xlsx -> rename sheet, deduplicate T.SEG/TSEGUI columns
#
# Commit: "Eliminada columna T.SEG y renombrada a TSEGUI"
#   - Sheet "plantilla_datos_pacientes" renamed to "Pacientes"
#   - Column AT ("T.SEG") is renamed to "TSEGUI"
#   - The old, separate "TSEGUI" column (AW) is deleted (duplicate header)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab.
$ws.Name = "Pacientes"

# Remove the duplicate "TSEGUI" header column; everything to its right
# shifts one column to the left (AX1:BM1 -> AW1:BL1).
$ws.Range("AW1").EntireColumn.Delete()

# Rename the original "T.SEG" header (now still at AT1) to "TSEGUI".
$ws.Range("AT1").Value = "TSEGUI"

# Mirror the saved selection/viewport from the edited workbook.
$ws.Application.Goto($ws.Range("AW6"))
